$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the current selection (diff drops the <selection> element, leaving a bare sheetView)
$ws.Range("A1").Select()

# Add the new row of data
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "dddd"
